$d = $word.ActiveDocument

function Get-ParaByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        # A paragraph's Range.Text always ends with its paragraph mark (CR);
        # strip that off before comparing against the expected plain text.
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1) Merge the two adjacent runs "Features" + "Information" (which sit
#    inside the [FeaturesInformation] placeholder paragraph) into a
#    single run "FeaturesInformation". Replacing the paragraph's text
#    with the identical text causes Word to coalesce the previously
#    split, identically-formatted runs into a single run.
# ------------------------------------------------------------------
$featuresPara = Get-ParaByText $d "[FeaturesInformation]"
if ($featuresPara -ne $null) {
    $null = $featuresPara.Range.Find.Execute("FeaturesInformation", $true, $false, $false, $false, $false, $true, 1, $false, "FeaturesInformation", 2)
}

# ------------------------------------------------------------------
# 2) Insert a brand new "[PublicEqualityDuty]" heading paragraph right
#    after the existing "[LegalInformation]" paragraph (i.e. directly
#    before the "[RationaleInformation]" paragraph).
# ------------------------------------------------------------------
$legalPara = Get-ParaByText $d "[LegalInformation]"
if ($legalPara -ne $null) {
    # Create a fresh empty paragraph immediately after [LegalInformation];
    # it inherits the Heading1 style used throughout this document.
    $insertionPoint = $legalPara.Range.End
    $collapsed = $d.Range($insertionPoint, $insertionPoint)
    $collapsed.InsertParagraphBefore()

    # The blank paragraph we just created is the one directly following
    # [LegalInformation]. ($legalPara itself stays valid/stable because
    # the insertion happened after its own end.)
    $newPara = $legalPara.Next()

    $fragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="24"/></w:rPr><w:t>[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="24"/></w:rPr><w:t>PublicEqualityDuty</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>]</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $null = $newPara.Range.InsertXML($fragment)
}
